# Add a new "Emner" worksheet (manual mapping of LDA topics) as the last sheet.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Emner"

# Column widths (characters), matching source column layout:
# A = EmneNr, B = Emne, C = Noter, D = SidstOpdatret
$ws.Columns.Item(1).ColumnWidth = 11.33203125
$ws.Columns.Item(2).ColumnWidth = 50.21875
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 21.33203125

# Header row (entered Emne, EmneNr, Noter, SidstOpdatret in that order)
$ws.Range("B1").Value = "Emne"
$ws.Range("A1").Value = "EmneNr"
$ws.Range("C1").Value = "Noter"
$ws.Range("D1").Value = "SidstOpdatret"

# Topic numbers (column A), entered top to bottom first
$ws.Range("A2").Value = "Topic 0"
$ws.Range("A3").Value = "Topic 1"
$ws.Range("A4").Value = "Topic 2"
$ws.Range("A5").Value = "Topic 3"
$ws.Range("A6").Value = "Topic 4"

# Topic descriptions / manual mapping (column B), entered after column A
$ws.Range("B2").Value = "Regeringen overholder ikke de eksisterende indfødsretsaftaler"
$ws.Range("B3").Value = 'Ansøgerne har arbejdet hårdt for at opfylde de strenge krav og fortjener en "tillykke"'
$ws.Range("B4").Value = "Der er forskellige tilføjelser/ændringsforslag under debatten"
$ws.Range("B5").Value = "Nogle mener, at Danmark burde bryde med den internationale konvention om statsløse"
$ws.Range("B6").Value = "Ingen flere korte bemærkninger inden afstemning"

# Last-updated date for every topic row
$ws.Range("D2").Value = 45566
$ws.Range("D3").Value = 45566
$ws.Range("D4").Value = 45566
$ws.Range("D5").Value = 45566
$ws.Range("D6").Value = 45566

# Reuse the date format already used on the Stopwords sheet (DatoTilføjet column)
$dateFormatSrc = $wb.Worksheets.Item("Stopwords").Range("B2")
$dateFormatSrc.Copy()
$ws.Range("D2:D6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the two-line topic descriptions
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8

# Alignment: column B wraps and is left/top aligned, the rest of the data
# rows are top-aligned.
$ws.Range("B2:B6").WrapText = $true
$ws.Range("B2:B6").HorizontalAlignment = -4131
$ws.Range("A2:D6").VerticalAlignment = -4160

# Turn the range into a proper table, like the other sheets in this workbook.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:D6"), $null, 1)
$lo.Name = "Emner"
$lo.TableStyle = "TableStyleMedium1"

# Page setup to match the other sheets.
$ws.PageSetup.Orientation = 1

# Make "Emner" the active sheet/tab with the same selection Excel left behind.
$ws.Activate()
$ws.Range("B12").Select()

Write-Output "Emner sheet added"
